# Auto-generated Excel COM-interop script applying updated market price data
# per the scheduled runner update (commit: chore: update Sheets via scheduled runner).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1015.6667
$ws.Range("I2").Value = 918.8
$ws.Range("K2").Value = 918.8
$ws.Range("M2").Value = -805.8
$ws.Range("H94").Value = 979.5
$ws.Range("I94").Value = 979.5
$ws.Range("K94").Value = 979.5
$ws.Range("M94").Value = -528.5
$ws.Range("H137").Value = 3973
$ws.Range("I137").Value = 3518.111
$ws.Range("J137").Value = 4996.5
$ws.Range("K137").Value = 10554.333
$ws.Range("L137").Value = 14989.5
$ws.Range("M137").Value = -8004.332999999999
$ws.Range("N137").Value = -20089.5
$ws.Range("H138").Value = 5779.25
$ws.Range("I138").Value = 3998.6667
$ws.Range("J138").Value = 6033.619
$ws.Range("K138").Value = 11996.0001
$ws.Range("L138").Value = 18100.857
$ws.Range("M138").Value = -6856.000100000001
$ws.Range("N138").Value = -28380.857

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4151.1875
$ws.Range("I32").Value = 3920.3618
$ws.Range("K32").Value = 3920.3618
$ws.Range("M32").Value = -3633.3618
$ws.Range("H74").Value = 5556068.5
$ws.Range("I74").Value = 6451370
$ws.Range("K74").Value = 6451370
$ws.Range("M74").Value = -6450496
$ws.Range("H77").Value = 5556068.5
$ws.Range("I77").Value = 6451370
$ws.Range("K77").Value = 32256850
$ws.Range("M77").Value = -32252482
$ws.Range("H110").Value = 1646.4
$ws.Range("I110").Value = 1837
$ws.Range("J110").Value = 1360.5
$ws.Range("K110").Value = 1837
$ws.Range("L110").Value = 1360.5
$ws.Range("M110").Value = 208
$ws.Range("N110").Value = -5450.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3307.5
$ws.Range("I134").Value = 3383.6956
$ws.Range("J134").Value = 1555
$ws.Range("K134").Value = 10151.0868
$ws.Range("L134").Value = 4665
$ws.Range("M134").Value = -7616.086800000001
$ws.Range("N134").Value = -9735
$ws.Range("H140").Value = 74999.5
$ws.Range("J140").Value = 75000
$ws.Range("L140").Value = 75000
$ws.Range("N140").Value = -85360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 755.5
$ws.Range("I16").Value = 755.5
$ws.Range("K16").Value = 755.5
$ws.Range("M16").Value = -468.5
$ws.Range("H31").Value = 2471.4473
$ws.Range("I31").Value = 2368.6572
$ws.Range("K31").Value = 2368.6572
$ws.Range("M31").Value = -2073.6572
$ws.Range("H34").Value = 2471.4473
$ws.Range("I34").Value = 2368.6572
$ws.Range("K34").Value = 2368.6572
$ws.Range("M34").Value = -2166.6572
$ws.Range("H62").Value = 2500
$ws.Range("I62").Value = 2500
$ws.Range("K62").Value = 2500
$ws.Range("M62").Value = -1876
$ws.Range("H65").Value = 2500
$ws.Range("I65").Value = 2500
$ws.Range("K65").Value = 12500
$ws.Range("M65").Value = -9380
$ws.Range("H68").Value = 49995
$ws.Range("J68").Value = 49995
$ws.Range("L68").Value = 49995
$ws.Range("N68").Value = -51493
$ws.Range("H71").Value = 49995
$ws.Range("J71").Value = 49995
$ws.Range("L71").Value = 149985
$ws.Range("N71").Value = -157473
$ws.Range("H74").Value = 55058.2
$ws.Range("J74").Value = 55058.2
$ws.Range("L74").Value = 55058.2
$ws.Range("N74").Value = -56806.2
$ws.Range("H77").Value = 55058.2
$ws.Range("J77").Value = 55058.2
$ws.Range("L77").Value = 165174.6
$ws.Range("N77").Value = -173910.6
$ws.Range("H113").Value = 755.5
$ws.Range("I113").Value = 755.5
$ws.Range("K113").Value = 755.5
$ws.Range("M113").Value = 1414.5
$ws.Range("H134").Value = 2215.8572
$ws.Range("I134").Value = 2215.8572
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 6647.571599999999
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -4112.571599999999
$ws.Range("H140").Value = 99999
$ws.Range("J140").Value = 99999
$ws.Range("L140").Value = 99999
$ws.Range("N140").Value = -110359

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 101180.1
$ws.Range("I122").Value = 904
$ws.Range("J122").Value = 112321.89
$ws.Range("K122").Value = 8136
$ws.Range("L122").Value = 1010897.01
$ws.Range("M122").Value = -5686
$ws.Range("N122").Value = -1015797.01
$ws.Range("H132").Value = 4783.4287
$ws.Range("I132").Value = 4663
$ws.Range("J132").Value = 4873.75
$ws.Range("K132").Value = 41967
$ws.Range("L132").Value = 43863.75
$ws.Range("M132").Value = -39437
$ws.Range("N132").Value = -48923.75
$ws.Range("H136").Value = 4075
$ws.Range("I136").Value = 4075
$ws.Range("K136").Value = 12225
$ws.Range("M136").Value = -7125
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").ClearContents()
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = 0

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 89.666664
$ws.Range("I2").Value = 82.125
$ws.Range("J2").Value = 150
$ws.Range("K2").Value = 82.125
$ws.Range("L2").Value = 150
$ws.Range("M2").Value = 30.875
$ws.Range("N2").Value = -376
$ws.Range("H80").Value = 5813
$ws.Range("I80").Value = 2921.5
$ws.Range("J80").Value = 9668.333000000001
$ws.Range("K80").Value = 2921.5
$ws.Range("L80").Value = 9668.333000000001
$ws.Range("M80").Value = -1923.5
$ws.Range("N80").Value = -11664.333
$ws.Range("H83").Value = 5813
$ws.Range("I83").Value = 2921.5
$ws.Range("J83").Value = 9668.333000000001
$ws.Range("K83").Value = 14607.5
$ws.Range("L83").Value = 48341.665
$ws.Range("M83").Value = -9615.5
$ws.Range("N83").Value = -58325.665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 166.66667
$ws.Range("J46").Value = 200
$ws.Range("L46").Value = 200
$ws.Range("N46").Value = -576
$ws.Range("H68").Value = 2951
$ws.Range("I68").Value = 2951
$ws.Range("K68").Value = 2951
$ws.Range("M68").Value = -2202
$ws.Range("H71").Value = 2951
$ws.Range("I71").Value = 2951
$ws.Range("K71").Value = 14755
$ws.Range("M71").Value = -11011
$ws.Range("H132").Value = 2767.6897
$ws.Range("I132").Value = 1798.7222
$ws.Range("J132").Value = 4353.273
$ws.Range("K132").Value = 5396.1666
$ws.Range("L132").Value = 13059.819
$ws.Range("M132").Value = -2866.1666
$ws.Range("N132").Value = -18119.819

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2011.9445
$ws.Range("I132").Value = 1428.079
$ws.Range("J132").Value = 3398.625
$ws.Range("K132").Value = 4284.237
$ws.Range("L132").Value = 10195.875
$ws.Range("M132").Value = -1754.237
$ws.Range("N132").Value = -15255.875
